# remove duplicate entries from bootstrap file
#
# Row 60 on the "hotels" sheet is an exact duplicate of row 59 (same
# placeid/fullname/displayname/address/locality/region/country/postal
# values and same lat/long). Delete the whole row and shift the rows
# below it up, which re-numbers the former rows 61-108 as 60-107.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("hotels")

$ws.Rows.Item(60).Delete()

# The two sheet-scoped defined names anchored on this sheet need to be
# shifted up by one row to keep tracking the same logical ranges
# (datafile: header block, hotels: data block).
$wb.Names.Item("hotels!datafile").RefersTo = "=hotels!`$A`$1:`$K`$66"
$wb.Names.Item("hotels!hotels").RefersTo = "=hotels!`$A`$67:`$K`$107"

# Reflect the scroll position / selection left behind by the edit.
$ws.Application.ActiveWindow.ScrollRow = 41
$ws.Range("A60:XFD60").Select()
